$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(44330,44331,44332,44333,44334,44335,44336,44337,44338,44339,44340,44341,44342,44343)
$bvals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$cvals = @(1,0,0,0,0,0,0,0,0,0,0,0,0,0)
$dvals = @(40.79967360261118,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 256 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
}

# Copy formatting (style) of the last original row (A255) down into the new A column cells
$ws.Range("A255").Copy()
$ws.Range("A256:A269").PasteSpecial(-4122)
$excel.CutCopyMode = $false
